$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update BF column (Date) cells: keep as text, not auto-converted to date serials ---
$ws.Range("BF2:BF31").NumberFormat = "@"

$ws.Range("BF2").Value = "2014-04-03"
$ws.Range("BF3").Value = "2014-04-03"
$ws.Range("BF4").Value = "2014-04-03"
$ws.Range("BF5").Value = "2014-04-03"
$ws.Range("BF6").Value = "2014-04-03"
$ws.Range("BF7").Value = "2014-04-03"
$ws.Range("BF8").Value = "2014-04-03"
$ws.Range("BF9").Value = "2014-04-03"
$ws.Range("BF10").Value = "2014-04-03"
$ws.Range("BF11").Value = "2014-04-03"
$ws.Range("BF12").Value = "2014-04-03"
$ws.Range("BF13").Value = "2014-04-03"
$ws.Range("BF14").Value = "2014-04-03"
$ws.Range("BF15").Value = "2014-04-03"
$ws.Range("BF16").Value = "2014-04-03"
$ws.Range("BF17").Value = "2014-04-03"
$ws.Range("BF18").Value = "2014-04-03"
$ws.Range("BF19").Value = "2014-04-03"
$ws.Range("BF20").Value = "2014-04-03"
$ws.Range("BF21").Value = "2014-04-03"
$ws.Range("BF22").Value = "2014-04-03"
$ws.Range("BF23").Value = "2014-04-03"
$ws.Range("BF24").Value = "2014-04-03"
$ws.Range("BF25").Value = "2014-04-03"
$ws.Range("BF26").Value = "2014-04-03"
$ws.Range("BF27").Value = "2014-04-03"
$ws.Range("BF28").Value = "2014-04-03"
$ws.Range("BF29").Value = "2014-04-03"
$ws.Range("BF30").Value = "2014-04-03"
$ws.Range("BF31").Value = "2014-04-03"

# --- Update numeric stat cells ---
$ws.Range("AN2").Value = 13
$ws.Range("AD3").Value = 6
$ws.Range("AX3").Value = 21
$ws.Range("AW4").Value = 6
$ws.Range("AD5").Value = 6
$ws.Range("AD6").Value = 6
$ws.Range("AX6").Value = 8
$ws.Range("AD7").Value = 1
$ws.Range("D8").Value = 75
$ws.Range("E8").Value = 44
$ws.Range("G8").Value = 0.587
$ws.Range("I8").Value = 39.5
$ws.Range("L8").Value = 8.699999999999999
$ws.Range("N8").Value = 0.384
$ws.Range("Q8").Value = 0.796
$ws.Range("R8").Value = 10.2
$ws.Range("S8").Value = 30.5
$ws.Range("T8").Value = 40.7
$ws.Range("AB8").Value = 105.1
$ws.Range("AD8").Value = 6
$ws.Range("AH8").Value = 18
$ws.Range("AK8").Value = 5
$ws.Range("AW8").Value = 3
$ws.Range("AD9").Value = 6
$ws.Range("AT9").Value = 6
$ws.Range("AV9").Value = 27
$ws.Range("AD10").Value = 6
$ws.Range("AV10").Value = 13
$ws.Range("AD11").Value = 6
$ws.Range("AH12").Value = 17
$ws.Range("AK12").Value = 4
$ws.Range("AD13").Value = 1
$ws.Range("AF13").Value = 5
$ws.Range("D14").Value = 76
$ws.Range("F14").Value = 22
$ws.Range("G14").Value = 0.711
$ws.Range("J14").Value = 82.3
$ws.Range("L14").Value = 8.5
$ws.Range("N14").Value = 0.354
$ws.Range("O14").Value = 21.1
$ws.Range("Q14").Value = 0.732
$ws.Range("V14").Value = 14
$ws.Range("Y14").Value = 3.4
$ws.Range("AC14").Value = 7.1
$ws.Range("AE14").Value = 2
$ws.Range("AF14").Value = 3
$ws.Range("AG14").Value = 3
$ws.Range("AL14").Value = 11
$ws.Range("AN14").Value = 18
$ws.Range("AV14").Value = 8
$ws.Range("AW14").Value = 5
$ws.Range("AD15").Value = 6
$ws.Range("AD16").Value = 6
$ws.Range("AE16").Value = 9
$ws.Range("AG16").Value = 9
$ws.Range("AN16").Value = 19
$ws.Range("AG17").Value = 4
$ws.Range("AD18").Value = 6
$ws.Range("AN18").Value = 20
$ws.Range("AW19").Value = 4
$ws.Range("AD20").Value = 6
$ws.Range("AV20").Value = 7
$ws.Range("AD21").Value = 1
$ws.Range("AP21").Value = 29
$ws.Range("D22").Value = 73
$ws.Range("E22").Value = 54
$ws.Range("G22").Value = 0.74
$ws.Range("M22").Value = 22.2
$ws.Range("Q22").Value = 0.804
$ws.Range("S22").Value = 34.3
$ws.Range("T22").Value = 45.1
$ws.Range("V22").Value = 15.8
$ws.Range("W22").Value = 8.300000000000001
$ws.Range("Z22").Value = 22.6
$ws.Range("AB22").Value = 106.2
$ws.Range("AC22").Value = 6.7
$ws.Range("AD22").Value = 30
$ws.Range("AL22").Value = 14
$ws.Range("AN22").Value = 11
$ws.Range("AT22").Value = 5
$ws.Range("AV22").Value = 28
$ws.Range("AW22").Value = 10
$ws.Range("AD23").Value = 6
$ws.Range("AD24").Value = 6
$ws.Range("AD25").Value = 6
$ws.Range("AE25").Value = 9
$ws.Range("AG25").Value = 9
$ws.Range("AW25").Value = 9
$ws.Range("AD26").Value = 1
$ws.Range("AD27").Value = 6
$ws.Range("D28").Value = 75
$ws.Range("F28").Value = 16
$ws.Range("G28").Value = 0.787
$ws.Range("J28").Value = 83
$ws.Range("K28").Value = 0.49
$ws.Range("N28").Value = 0.399
$ws.Range("O28").Value = 15.9
$ws.Range("P28").Value = 20.3
$ws.Range("U28").Value = 25.4
$ws.Range("W28").Value = 7.5
$ws.Range("AB28").Value = 105.6
$ws.Range("AC28").Value = 8.6
$ws.Range("AD28").Value = 6
$ws.Range("AH28").Value = 24
$ws.Range("AL28").Value = 12
$ws.Range("AV28").Value = 12
$ws.Range("AX28").Value = 7
$ws.Range("AD29").Value = 6
$ws.Range("AN29").Value = 12
$ws.Range("AD30").Value = 6
$ws.Range("AD31").Value = 6
$ws.Range("AL31").Value = 15
